# Fruta / hortaliza, semanal
#
# Inserts a new weekly price-record row for "Vega Modelo de Temuco - Ciboulette"
# at sheet row 187 (pushing the existing rows 187:229 down to 188:230), then
# fills the new row with the same record values as the row that used to sit
# at 187, except for the date (column D) and the volume (column J), which get
# the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 187; everything below shifts down
# by one (old row 187 becomes row 188, old 229 becomes 230, dimension grows
# from A1:R229 to A1:R230 automatically).
$ws.Rows(187).Insert()

# The record that used to live in row 187 is now in row 188 - duplicate it
# into the freshly inserted row 187 as the starting point for the new entry.
$ws.Range("A187:R187").Value = $ws.Range("A188:R188").Value()

# Overwrite the two fields that differ for this new weekly record.
$ws.Range("D187").Value = 44642
$ws.Range("J187").Value = 65
